# Update countries & provincias Spain
#
# Applies the 10-Abril-2020 data refresh:
#  - Re-sorts three country pairs/groups whose case counts changed enough to
#    change their rank in the (descending) table: Ecuador/Japon,
#    Guatemala/Paraguay, and Birmania (which leap-frogs Benin/Tanzania/Libia).
#  - Refreshes the numeric stats (Casos totales, Nuevos casos, Casos activos,
#    Recuperados, Casos criticos, Muertes hoy, Muertes) for the affected rows
#    plus a couple of standalone updates (Australia, Bermudas).
#  - Bumps the "Datos actualizados" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($Row, $Pais, $CasosTotales, $NuevosCasos, $CasosActivos, $Recuperados, $CasosCriticos, $MuertesHoy, $Muertes) {
    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $CasosTotales
    $ws.Cells.Item($Row, 3).Value = $NuevosCasos
    $ws.Cells.Item($Row, 4).Value = $CasosActivos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $CasosCriticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# --- Ecuador / Japon swap places (row 34 <-> 35), both refreshed with today's
#     counts for Japon (now ranked above Ecuador) ---
Set-CountryRow 34 "Japon"   4979 0  685 4195 109 0 99
Set-CountryRow 35 "Ecuador" 4965 0  339 4354 139 0 272

# --- Guatemala / Paraguay swap places (row 125 <-> 126); Paraguay refreshed ---
Set-CountryRow 125 "Paraguay"  129 5  18 105 1 1 6
Set-CountryRow 126 "Guatemala" 126 31 17 106 3 0 3

# --- Birmania jumps ahead of Benin/Tanzania/Libia (rows 160-163), each of the
#     latter three simply shifts down one row; Birmania is refreshed ---
Set-CountryRow 160 "Birmania" 27 4 2 22 0 0 3
Set-CountryRow 161 "Benin"    26 0 5 20 0 0 1
Set-CountryRow 162 "Tanzania" 25 0 5 19 0 0 1
Set-CountryRow 163 "Libia"    24 0 8 15 0 0 1

# --- Standalone numeric refreshes (no reordering) ---
# Australia (row 27): Recuperados, Muertes hoy, Muertes
$ws.Cells.Item(27, 5).Value = 3113
$ws.Cells.Item(27, 7).Value = 1
$ws.Cells.Item(27, 8).Value = 52

# Bermudas (row 145): Casos criticos
$ws.Cells.Item(145, 6).Value = 2

# --- Footer timestamp bump ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 03:22"
